# Natmi following Dr Hou advice
# Rebuild the LR-pair table body (rows 2-10) on the active sheet with the
# recomputed sending-cluster / target-cluster combinations (ECs, FAPs, sCs)
# and their updated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Wnt11", "Fzd7", "ECs", 1, 0.3333333333333333, 0.1470843333333333, 0.441253, 0.01080519019772543, 0.01080519019772543, 3, 1, 2.553279333333334, 7.659838000000001, 0.1645043904057808, 0.1645043904057808, 0.3755473885571112, 3.379926497014, 0.001777501226695341, 0.00177750122669534),
    @("ECs", "Wnt11", "Fzd7", "FAPs", 1, 0.3333333333333333, 0.1470843333333333, 0.441253, 0.01080519019772543, 0.01080519019772543, 3, 1, 8.058662, 24.175986, 0.5192088709172035, 0.5192088709172035, 1.185302927828667, 10.667726350458, 0.005610150602606657, 0.005610150602606657),
    @("ECs", "Wnt11", "Fzd7", "sCs", 1, 0.3333333333333333, 0.1470843333333333, 0.441253, 0.01080519019772543, 0.01080519019772543, 3, 1, 4.909099333333334, 14.727298, 0.3162867386770157, 0.3162867386770157, 0.7220516027104446, 6.498464424394, 0.003417538368423436, 0.003417538368423436),
    @("FAPs", "Wnt11", "Fzd7", "ECs", 3, 1, 13.22273466666667, 39.668204, 0.9713758071269154, 0.9713758071269153, 3, 1, 2.553279333333334, 7.659838000000001, 0.1645043904057808, 0.1645043904057808, 33.76133515455022, 303.852016390952, 0.1597955850063365, 0.1597955850063365),
    @("FAPs", "Wnt11", "Fzd7", "FAPs", 3, 1, 13.22273466666667, 39.668204, 0.9713758071269154, 0.9713758071269153, 3, 1, 8.058662, 24.175986, 0.5192088709172035, 0.5192088709172035, 106.5575493943493, 959.017944549144, 0.504346936054653, 0.5043469360546529),
    @("FAPs", "Wnt11", "Fzd7", "sCs", 3, 1, 13.22273466666667, 39.668204, 0.9713758071269154, 0.9713758071269153, 3, 1, 4.909099333333334, 14.727298, 0.3162867386770157, 0.3162867386770157, 64.91171793697688, 584.2054614327919, 0.3072332860659259, 0.3072332860659258),
    @("sCs", "Wnt11", "Fzd7", "ECs", 3, 1, 0.242559, 0.7276769999999999, 0.01781900267535914, 0.01781900267535914, 3, 1, 2.553279333333334, 7.659838000000001, 0.1645043904057808, 0.1645043904057808, 0.619320881814, 5.573887936326, 0.002931304172748933, 0.002931304172748933),
    @("sCs", "Wnt11", "Fzd7", "FAPs", 3, 1, 0.242559, 0.7276769999999999, 0.01781900267535914, 0.01781900267535914, 3, 1, 8.058662, 24.175986, 0.5192088709172035, 0.5192088709172035, 1.954700996058, 17.592308964522, 0.009251784259943848, 0.009251784259943848),
    @("sCs", "Wnt11", "Fzd7", "sCs", 3, 1, 0.242559, 0.7276769999999999, 0.01781900267535914, 0.01781900267535914, 3, 1, 4.909099333333334, 14.727298, 0.3162867386770157, 0.3162867386770157, 1.190746225194, 10.716716026746, 0.005635914242666361, 0.005635914242666361)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $excelRow = $i + 2
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowData[$c]
    }
}
